$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order IDs updated) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16511686782573843"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168681466479"
$wb.Worksheets.Item(3).Name = "RS_TO-16511686814674158"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511686815154202"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511686815944245"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511686782215765.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686782404156.csv"
$ws1.Range("B4").Value = "go_stims-16511686782413797.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686782564163.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16511686805014184.csv"
$ws2.Range("B3").Value = "OB-16511686793063512.csv"
$ws2.Range("B4").Value = "TB-16511686814504569.csv"
$ws2.Range("B5").Value = "OB-16511686790303876.csv"
$ws2.Range("B6").Value = "TB-16511686809174204.csv"
$ws2.Range("B7").Value = "ZB-match_1-16511686786173558.csv"
$ws2.Range("B8").Value = "OB-1651168680156388.csv"
$ws2.Range("B9").Value = "ZB-match_0-16511686785543888.csv"
$ws2.Range("B10").Value = "ZB-match_0-1651168678710387.csv"

# --- Sheet 3 (RS) unchanged ---

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511686814824507.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686814704287.csv"
$ws4.Range("B4").Value = "MM_stims-1651168681498453.csv"
$ws4.Range("B5").Value = "ZM_stims-165116868148342.csv"
$ws4.Range("B6").Value = "MM_stims-16511686815144572.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686814994278.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1651168681578419.csv"
$ws5.Range("B3").Value = "SAT_stims-1651168681546421.csv"
$ws5.Range("B4").Value = "SAT_stims-16511686815204191.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511686815624557.csv"
